$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Actual Outcome" (col F) and "Fail/Pass" (col G) columns for
# test case rows 2-4: the tests have now been run, and all three passed
# (previously these cells just held placeholder "-" / "Fail" values).
$ws.Range("F2:F4").Value = "Same as expected outcome."
$ws.Range("G2:G4").Value = "Pass"

# Reflect the updated selection/scroll position left by the author after
# reviewing the test outcomes: active cell H4, scrolled so column D is the
# left-most visible column.
$ws.Activate()
$ws.Range("H4").Select()
$excel.ActiveWindow.ScrollColumn = 4

$wb.Save()
